$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.011373231998788
$ws.Cells.Item(2, 4).Value = 1.032035228338414
$ws.Cells.Item(2, 5).Value = 1.013534140562208
$ws.Cells.Item(2, 6).Value = 1.025684107681387
$ws.Cells.Item(2, 9).Value = 1.031307747648028
$ws.Cells.Item(2, 10).Value = 1.016622080019619
$ws.Cells.Item(2, 11).Value = 1.034841621578659
$ws.Cells.Item(2, 12).Value = 1.016395004300311
$ws.Cells.Item(2, 13).Value = 1.028508958798792
$ws.Cells.Item(2, 14).Value = 1.01806579998143

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.01258891769867
$ws.Cells.Item(3, 4).Value = 1.032452732606883
$ws.Cells.Item(3, 5).Value = 1.014573472503753
$ws.Cells.Item(3, 6).Value = 1.027029552794126
$ws.Cells.Item(3, 9).Value = 1.031373466714418
$ws.Cells.Item(3, 10).Value = 1.017469576625707
$ws.Cells.Item(3, 11).Value = 1.035068862086806
$ws.Cells.Item(3, 12).Value = 1.017238433326009
$ws.Cells.Item(3, 13).Value = 1.029660299467545
$ws.Cells.Item(3, 14).Value = 1.018914500129908

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.01337547087054
$ws.Cells.Item(4, 4).Value = 1.032722402689852
$ws.Cells.Item(4, 5).Value = 1.01524624216872
$ws.Cells.Item(4, 6).Value = 1.027899267877749
$ws.Cells.Item(4, 9).Value = 1.031414256784638
$ws.Cells.Item(4, 10).Value = 1.01801744984403
$ws.Cells.Item(4, 11).Value = 1.035214684451568
$ws.Cells.Item(4, 12).Value = 1.017783849402546
$ws.Cells.Item(4, 13).Value = 1.030403892602807
$ws.Cells.Item(4, 14).Value = 1.019463151391042

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.013706122056136
$ws.Cells.Item(5, 4).Value = 1.032835654442424
$ws.Cells.Item(5, 5).Value = 1.015529136345781
$ws.Cells.Item(5, 6).Value = 1.02826468856932
$ws.Cells.Item(5, 9).Value = 1.031430989025182
$ws.Cells.Item(5, 10).Value = 1.018247654262254
$ws.Cells.Item(5, 11).Value = 1.035275695452555
$ws.Cells.Item(5, 12).Value = 1.018013062514757
$ws.Cells.Item(5, 13).Value = 1.030716165782586
$ws.Cells.Item(5, 14).Value = 1.019693682725946

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.013761638996323
$ws.Cells.Item(6, 4).Value = 1.032854662961379
$ws.Cells.Item(6, 5).Value = 1.015576639187636
$ws.Cells.Item(6, 6).Value = 1.028326032218047
$ws.Cells.Item(6, 9).Value = 1.031433774034169
$ws.Cells.Item(6, 10).Value = 1.018286299517072
$ws.Cells.Item(6, 11).Value = 1.035285922265813
$ws.Cells.Item(6, 12).Value = 1.018051543763504
$ws.Cells.Item(6, 13).Value = 1.030768578334837
$ws.Cells.Item(6, 14).Value = 1.019732382861457

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.013379889106574
$ws.Cells.Item(7, 4).Value = 1.032723916429626
$ws.Cells.Item(7, 5).Value = 1.015250021969976
$ws.Cells.Item(7, 6).Value = 1.02790415145916
$ws.Cells.Item(7, 9).Value = 1.031414481996706
$ws.Cells.Item(7, 10).Value = 1.018020526322211
$ws.Cells.Item(7, 11).Value = 1.035215500835358
$ws.Cells.Item(7, 12).Value = 1.017786912471316
$ws.Cells.Item(7, 13).Value = 1.030408066519731
$ws.Cells.Item(7, 14).Value = 1.019466232238175

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.011784096139778
$ws.Cells.Item(8, 4).Value = 1.032176424753418
$ws.Cells.Item(8, 5).Value = 1.013885336231729
$ws.Cells.Item(8, 6).Value = 1.026138990515252
$ws.Cells.Item(8, 9).Value = 1.031330316478876
$ws.Cells.Item(8, 10).Value = 1.016908602753379
$ws.Cells.Item(8, 11).Value = 1.034918669894545
$ws.Cells.Item(8, 12).Value = 1.016680116060586
$ws.Cells.Item(8, 13).Value = 1.028898351172271
$ws.Cells.Item(8, 14).Value = 1.018352729610336

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.00897140005891
$ws.Cells.Item(9, 4).Value = 1.031208057710918
$ws.Cells.Item(9, 5).Value = 1.011482456214916
$ws.Cells.Item(9, 6).Value = 1.023021709170121
$ws.Cells.Item(9, 9).Value = 1.031168749680305
$ws.Cells.Item(9, 10).Value = 1.014945248553917
$ws.Cells.Item(9, 11).Value = 1.034386347070717
$ws.Cells.Item(9, 12).Value = 1.014727143959245
$ws.Cells.Item(9, 13).Value = 1.026227213278787
$ws.Cells.Item(9, 14).Value = 1.016386587222711

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.007095629892131
$ws.Cells.Item(10, 4).Value = 1.030560159682543
$ws.Cells.Item(10, 5).Value = 1.009881704194111
$ws.Cells.Item(10, 6).Value = 1.020938751910089
$ws.Cells.Item(10, 9).Value = 1.03105216305053
$ws.Cells.Item(10, 10).Value = 1.013633557449773
$ws.Cells.Item(10, 11).Value = 1.034025307423701
$ws.Cells.Item(10, 12).Value = 1.013423291561006
$ws.Cells.Item(10, 13).Value = 1.024439039966655
$ws.Cells.Item(10, 14).Value = 1.015073033366746

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.006283211648611
$ws.Cells.Item(11, 4).Value = 1.030279085851953
$ws.Cells.Item(11, 5).Value = 1.009188817235273
$ws.Cells.Item(11, 6).Value = 1.020035638405816
$ws.Cells.Item(11, 9).Value = 1.030999582288564
$ws.Cells.Item(11, 10).Value = 1.013064897179479
$ws.Cells.Item(11, 11).Value = 1.033867528167864
$ws.Cells.Item(11, 12).Value = 1.012858245621728
$ws.Cells.Item(11, 13).Value = 1.023662951298119
$ws.Cells.Item(11, 14).Value = 1.014503565533642

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.005981410386023
$ws.Cells.Item(12, 4).Value = 1.030174604816601
$ws.Cells.Item(12, 5).Value = 1.008931483744683
$ws.Cells.Item(12, 6).Value = 1.019700001205032
$ws.Cells.Item(12, 9).Value = 1.030979736833185
$ws.Cells.Item(12, 10).Value = 1.012853565875572
$ws.Cells.Item(12, 11).Value = 1.03380870579526
$ws.Cells.Item(12, 12).Value = 1.012648290128377
$ws.Cells.Item(12, 13).Value = 1.023374404485686
$ws.Cells.Item(12, 14).Value = 1.014291934115043

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.00604614932566
$ws.Cells.Item(13, 4).Value = 1.030197019839373
$ws.Cells.Item(13, 5).Value = 1.008986681072115
$ws.Cells.Item(13, 6).Value = 1.019772004753836
$ws.Cells.Item(13, 9).Value = 1.030984007978691
$ws.Cells.Item(13, 10).Value = 1.012898901966501
$ws.Cells.Item(13, 11).Value = 1.033821333154826
$ws.Cells.Item(13, 12).Value = 1.01269332959348
$ws.Cells.Item(13, 13).Value = 1.023436311152422
$ws.Cells.Item(13, 14).Value = 1.014337334588421

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.00625826535347
$ws.Cells.Item(14, 4).Value = 1.030270450991202
$ws.Cells.Item(14, 5).Value = 1.009167545253745
$ws.Cells.Item(14, 6).Value = 1.020007898235087
$ws.Cells.Item(14, 9).Value = 1.030997948271846
$ws.Cells.Item(14, 10).Value = 1.013047430632908
$ws.Cells.Item(14, 11).Value = 1.033862670294211
$ws.Cells.Item(14, 12).Value = 1.012840892119786
$ws.Cells.Item(14, 13).Value = 1.02363910551259
$ws.Cells.Item(14, 14).Value = 1.014486074182571

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.00638895257332
$ws.Cells.Item(15, 4).Value = 1.030315684107427
$ws.Cells.Item(15, 5).Value = 1.009278986303561
$ws.Cells.Item(15, 6).Value = 1.020153215954617
$ws.Cells.Item(15, 9).Value = 1.031006495673641
$ws.Cells.Item(15, 10).Value = 1.013138929998941
$ws.Cells.Item(15, 11).Value = 1.033888110866782
$ws.Cells.Item(15, 12).Value = 1.012931800632363
$ws.Cells.Item(15, 13).Value = 1.023764017560877
$ws.Cells.Item(15, 14).Value = 1.014577703488199

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.007149542975283
$ws.Cells.Item(16, 4).Value = 1.030578802627726
$ws.Cells.Item(16, 5).Value = 1.009927693911421
$ws.Cells.Item(16, 6).Value = 1.020998663438169
$ws.Cells.Item(16, 9).Value = 1.031055608510028
$ws.Cells.Item(16, 10).Value = 1.013671282865223
$ws.Cells.Item(16, 11).Value = 1.034035748337195
$ws.Cells.Item(16, 12).Value = 1.013460781761361
$ws.Cells.Item(16, 13).Value = 1.024490508297005
$ws.Cells.Item(16, 14).Value = 1.015110812356613

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.007626586068018
$ws.Cells.Item(17, 4).Value = 1.030743709440082
$ws.Cells.Item(17, 5).Value = 1.010334676465269
$ws.Cells.Item(17, 6).Value = 1.021528672146503
$ws.Cells.Item(17, 9).Value = 1.031085854416069
$ws.Cells.Item(17, 10).Value = 1.014005027787794
$ws.Cells.Item(17, 11).Value = 1.03412797096188
$ws.Cells.Item(17, 12).Value = 1.013792470767327
$ws.Cells.Item(17, 13).Value = 1.024945733426109
$ws.Cells.Item(17, 14).Value = 1.015445031235254

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.007904818412832
$ws.Cells.Item(18, 4).Value = 1.030839845559506
$ws.Cells.Item(18, 5).Value = 1.01057208667803
$ws.Cells.Item(18, 6).Value = 1.021837703713386
$ws.Cells.Item(18, 9).Value = 1.031103293878347
$ws.Cells.Item(18, 10).Value = 1.014199629176507
$ws.Cells.Item(18, 11).Value = 1.034181623191605
$ws.Cells.Item(18, 12).Value = 1.013985894198398
$ws.Cells.Item(18, 13).Value = 1.025211085271136
$ws.Cells.Item(18, 14).Value = 1.015639908980261

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.007999685359256
$ws.Cells.Item(19, 4).Value = 1.030872616743965
$ws.Cells.Item(19, 5).Value = 1.010653041581723
$ws.Cells.Item(19, 6).Value = 1.021943056336174
$ws.Cells.Item(19, 9).Value = 1.031109205925251
$ws.Cells.Item(19, 10).Value = 1.014265972021963
$ws.Cells.Item(19, 11).Value = 1.034199893488104
$ws.Cells.Item(19, 12).Value = 1.014051838984098
$ws.Cells.Item(19, 13).Value = 1.025301534087788
$ws.Cells.Item(19, 14).Value = 1.015706346040168

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.007575405830645
$ws.Cells.Item(20, 4).Value = 1.030726021785794
$ws.Cells.Item(20, 5).Value = 1.01029100860286
$ws.Cells.Item(20, 6).Value = 1.021471818999161
$ws.Cells.Item(20, 9).Value = 1.031082630254263
$ws.Cells.Item(20, 10).Value = 1.013969226997228
$ws.Cells.Item(20, 11).Value = 1.034118090784824
$ws.Cells.Item(20, 12).Value = 1.013756888360281
$ws.Cells.Item(20, 13).Value = 1.024896910021958
$ws.Cells.Item(20, 14).Value = 1.01540917960346

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.006195803406372
$ws.Cells.Item(21, 4).Value = 1.030248829477068
$ws.Cells.Item(21, 5).Value = 1.009114284293127
$ws.Cells.Item(21, 6).Value = 1.019938438498648
$ws.Cells.Item(21, 9).Value = 1.03099385188338
$ws.Cells.Item(21, 10).Value = 1.013003695572629
$ws.Cells.Item(21, 11).Value = 1.033850503493526
$ws.Cells.Item(21, 12).Value = 1.012797440643059
$ws.Cells.Item(21, 13).Value = 1.023579395184883
$ws.Cells.Item(21, 14).Value = 1.01444227701349

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.005328200303288
$ws.Cells.Item(22, 4).Value = 1.029948351105582
$ws.Cells.Item(22, 5).Value = 1.008374635920009
$ws.Cells.Item(22, 6).Value = 1.018973294494842
$ws.Cells.Item(22, 9).Value = 1.030936213389767
$ws.Cells.Item(22, 10).Value = 1.012396015915169
$ws.Cells.Item(22, 11).Value = 1.033681010865718
$ws.Cells.Item(22, 12).Value = 1.012193778750151
$ws.Cells.Item(22, 13).Value = 1.022749440839736
$ws.Cells.Item(22, 14).Value = 1.013833734381215

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.005788152456636
$ws.Cells.Item(23, 4).Value = 1.030107682275835
$ws.Cells.Item(23, 5).Value = 1.008766718720886
$ws.Cells.Item(23, 6).Value = 1.019485035981519
$ws.Cells.Item(23, 9).Value = 1.030966940970075
$ws.Cells.Item(23, 10).Value = 1.012718217051296
$ws.Cells.Item(23, 11).Value = 1.033770980184498
$ws.Cells.Item(23, 12).Value = 1.012513831655473
$ws.Cells.Item(23, 13).Value = 1.023189566100286
$ws.Cells.Item(23, 14).Value = 1.014156393079912

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.007598532034943
$ws.Cells.Item(24, 4).Value = 1.030734014234042
$ws.Cells.Item(24, 5).Value = 1.010310740155553
$ws.Cells.Item(24, 6).Value = 1.021497508839954
$ws.Cells.Item(24, 9).Value = 1.031084087739895
$ws.Cells.Item(24, 10).Value = 1.013985404038529
$ws.Cells.Item(24, 11).Value = 1.034122555643315
$ws.Cells.Item(24, 12).Value = 1.013772966657916
$ws.Cells.Item(24, 13).Value = 1.024918971752018
$ws.Cells.Item(24, 14).Value = 1.015425379618015

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.009698649969702
$ws.Cells.Item(25, 4).Value = 1.031458820838607
$ws.Cells.Item(25, 5).Value = 1.012103445825459
$ws.Cells.Item(25, 6).Value = 1.023828429331001
$ws.Cells.Item(25, 9).Value = 1.031212085640675
$ws.Cells.Item(25, 10).Value = 1.015453306994268
$ws.Cells.Item(25, 11).Value = 1.034525056006715
$ws.Cells.Item(25, 12).Value = 1.015232357532454
$ws.Cells.Item(25, 13).Value = 1.026919063347634
$ws.Cells.Item(25, 14).Value = 1.016895367164322
